# Apply "development in TES class file" update to the TES Code Error Log
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FF Code Base")

# New column header for existing data
$ws.Range("E1").Value = "Relevant Code"
$ws.Range("E2").Value = "FF_test.m"

# New log entry (row 3)
$ws.Range("A3").Value = 44795
$ws.Range("A3").NumberFormat = "m/d/yyyy"
$ws.Range("B3").Value = "Instability durring simulation, originating from oscilations at the wall and base"
$ws.Range("C3").Value = "The scaling associated with transfering data from wall/base to the particle domain was incorrect. Initial scaling should use prototype parameters and the scaling used to transfer information to particle domain should use the model parameters."
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = 44795
$ws.Range("D3").NumberFormat = "m/d/yyyy"
$ws.Range("E3").Value = "TES.m, FF.m"

$ws.Rows.Item(3).RowHeight = 47.25

# Column width for the new column
$ws.Columns.Item(5).ColumnWidth = 16.875

# Update selection to match the saved view state
$ws.Range("E7").Select()
